{"js": "// Apply the \"storage html format made\" edit to the MY STORAGE paragraph.\n// Three targeted, unique-in-document search/replace operations that together\n// reproduce the run-level diff:\n//   1) \"This is a place ... from Add Food page.\" -> \"This is ... from the Add Food page.\"\n//   2) \"There is a search bar ... take it from the list.\" -> \"There is a sort bar ... easily update the list.\"\n//   3) \" out. One huge key ... My Storage page. \" -> \" out. \"  (trailing icon/logo explanation removed)\n\nconst replacements = [\n  {\n    find: \"This is a place where you can see all of the food you added from Add Food page.\",\n    replace: \"This is where you can see all of the food you added from the Add Food page.\"\n  },\n  {\n    find: \"There is a search bar to make it simpler, so that when an item is used up you can take it from the list.\",\n    replace: \"There is a sort bar to help with organization, so that when an item is used up etc. you can easily update the list.\"\n  },\n  {\n    find: \" out. One huge key that helps this website to make it so simple to use is the icons in the top left corner. The logo will take you to the home page, the plus logo to the Add Food page, and the Pantry icon to My Storage page. \",\n    replace: \" out. \"\n  }\n];\n\nfor (const { find, replace } of replacements) {\n  const results = context.document.body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Search text not found: \" + find);\n  }\n\n  results.items[0].insertText(replace, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Apply the \"storage html format made\" edit to the MY STORAGE paragraph.\n# Three targeted, unique Find/Replace operations scoped to that paragraph that\n# together reproduce the run-level diff:\n#   1) \"This is a place ... from Add Food page.\" -> \"This is ... from the Add Food page.\"\n#   2) \"There is a search bar ... take it from the list.\" -> \"There is a sort bar ... easily update the list.\"\n#   3) \" out. One huge key ... My Storage page. \" -> \" out. \"  (trailing icon/logo explanation removed)\n\n$d = $word.ActiveDocument\n\n# Locate the \"MY STORAGE\" body paragraph (the one starting with \"This is a place\")\n# and remember its index so we can keep re-fetching a fresh Range before each\n# Find/Replace call (Find.Execute collapses/advances the Range it runs on).\n$targetIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text.StartsWith(\"This is a place where you can see\")) {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -eq -1) {\n    throw \"Could not locate the 'This is a place where you can see ...' paragraph\"\n}\n\n$wdFindContinue = 1\n$wdReplaceOne = 1\n\n$range1 = $d.Paragraphs.Item($targetIndex).Range\n$find1 = $range1.Find\n$find1.Execute(\n    \"This is a place where you can see all of the food you added from Add Food page.\",\n    $false, $false, $false, $false, $false, $true, $wdFindContinue, $false,\n    \"This is where you can see all of the food you added from the Add Food page.\",\n    $wdReplaceOne\n) | Out-Null\n\n$range2 = $d.Paragraphs.Item($targetIndex).Range\n$find2 = $range2.Find\n$find2.Execute(\n    \"There is a search bar to make it simpler, so that when an item is used up you can take it from the list.\",\n    $false, $false, $false, $false, $false, $true, $wdFindContinue, $false,\n    \"There is a sort bar to help with organization, so that when an item is used up etc. you can easily update the list.\",\n    $wdReplaceOne\n) | Out-Null\n\n$range3 = $d.Paragraphs.Item($targetIndex).Range\n$find3 = $range3.Find\n$find3.Execute(\n    \" out. One huge key that helps this website to make it so simple to use is the icons in the top left corner. The logo will take you to the home page, the plus logo to the Add Food page, and the Pantry icon to My Storage page. \",\n    $false, $false, $false, $false, $false, $true, $wdFindContinue, $false,\n    \" out. \",\n    $wdReplaceOne\n) | Out-Null\n"}
